$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1
$ws.Cells.Item(1,6).Value = "time_taken"

# Copy style (bold, border, center/top alignment) from E1 header to F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate time_taken values for data rows 2-32
$ws.Cells.Item(2,6).Value = "2021-10-05 10:51:37.832674"
$ws.Cells.Item(3,6).Value = "2021-10-05 10:51:37.832687"
$ws.Cells.Item(4,6).Value = "2021-10-05 10:51:37.832691"
$ws.Cells.Item(5,6).Value = "2021-10-05 10:51:37.832695"
$ws.Cells.Item(6,6).Value = "2021-10-05 10:51:37.832698"
$ws.Cells.Item(7,6).Value = "2021-10-05 10:51:37.832701"
$ws.Cells.Item(8,6).Value = "2021-10-05 10:51:37.832704"
$ws.Cells.Item(9,6).Value = "2021-10-05 10:51:37.832707"
$ws.Cells.Item(10,6).Value = "2021-10-05 10:51:37.832711"
$ws.Cells.Item(11,6).Value = "2021-10-05 10:51:37.832714"
$ws.Cells.Item(12,6).Value = "2021-10-05 10:51:37.832717"
$ws.Cells.Item(13,6).Value = "2021-10-05 10:51:37.832720"
$ws.Cells.Item(14,6).Value = "2021-10-05 10:51:37.832723"
$ws.Cells.Item(15,6).Value = "2021-10-05 10:51:37.832725"
$ws.Cells.Item(16,6).Value = "2021-10-05 10:51:37.832728"
$ws.Cells.Item(17,6).Value = "2021-10-05 10:51:37.832731"
$ws.Cells.Item(18,6).Value = "2021-10-05 10:51:37.832735"
$ws.Cells.Item(19,6).Value = "2021-10-05 10:51:37.832738"
$ws.Cells.Item(20,6).Value = "2021-10-05 10:51:37.832741"
$ws.Cells.Item(21,6).Value = "2021-10-05 10:51:37.832744"
$ws.Cells.Item(22,6).Value = "2021-10-05 10:51:37.832748"
$ws.Cells.Item(23,6).Value = "2021-10-05 10:51:37.832750"
$ws.Cells.Item(24,6).Value = "2021-10-05 10:51:37.832753"
$ws.Cells.Item(25,6).Value = "2021-10-05 10:51:37.832756"
$ws.Cells.Item(26,6).Value = "2021-10-05 10:51:37.832760"
$ws.Cells.Item(27,6).Value = "2021-10-05 10:51:37.832763"
$ws.Cells.Item(28,6).Value = "2021-10-05 10:51:37.832766"
$ws.Cells.Item(29,6).Value = "2021-10-05 10:51:37.832769"
$ws.Cells.Item(30,6).Value = "2021-10-05 10:51:37.832772"
$ws.Cells.Item(31,6).Value = "2021-10-05 10:51:37.832775"
$ws.Cells.Item(32,6).Value = "2021-10-05 10:51:37.832778"

$excel.CutCopyMode = 0
